# fix up sentiment analysis
# Collapse the 5-category sentiment summary (Very Negative / Neutral / Positive /
# Very Positive / Negative) down to a simple 2-category Positive / Negative
# breakdown with updated counts, and drop the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "Positive" / "Very Positive" / "Negative" rows (4-6); rows 1-3
# are kept and overwritten below with the new, smaller summary.
$ws.Rows("4:6").Delete()

# Row 3 used to hold a labeled category ("Positive", count 9). The new layout
# only needs the total in column B on row 3, so fully clear A3 (value +
# formatting) rather than leaving an empty styled cell behind.
$ws.Range("A3").Clear()

# New 2-row sentiment labels.
$ws.Range("A1").Value = "Positive"
$ws.Range("A2").Value = "Negative"

# Updated counts: 1 positive, 1 negative, 98 unclassified/total in row 3.
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 98

# Reset the active selection back to the top of the sheet.
$ws.Range("A1").Select()
